# This script inserts three new data rows (for the "Sin especificar"
# quality entries dated serial 44498, i.e. 2021-10-29) into the
# "Hortaliza, Vega Central Mapocho de Santiago - Espárragos" weekly
# data sheet. All the existing data rows from 24 downward are pushed
# down by three rows (Excel keeps their content and formatting intact
# automatically), and the three freshly inserted rows (24-26) are then
# populated with the new observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows before row 24; rows 24-79 become 27-82.
$ws.Range("A24:A26").EntireRow.Insert()

# --- New row 24 ---
$ws.Cells.Item(24, 1).Value2 = 9
$ws.Cells.Item(24, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(24, 3).Value2 = "Metropolitana"
$ws.Cells.Item(24, 4).Value2 = 44498
$ws.Cells.Item(24, 5).Value2 = 13
$ws.Cells.Item(24, 6).Value2 = 300000000
$ws.Cells.Item(24, 7).Value2 = "Espárragos"
$ws.Cells.Item(24, 8).Value2 = "Sin especificar"
$ws.Cells.Item(24, 9).Value2 = "Banquete"
$ws.Cells.Item(24, 10).Value2 = 250
$ws.Cells.Item(24, 11).Value2 = 1200
$ws.Cells.Item(24, 12).Value2 = 1300
$ws.Cells.Item(24, 13).Value2 = 1240
$ws.Cells.Item(24, 14).Value2 = "`$/kilo"
$ws.Cells.Item(24, 15).Value2 = "Provincia de Linares"
$ws.Cells.Item(24, 16).Value2 = 1240
$ws.Cells.Item(24, 17).Value2 = 1
$ws.Cells.Item(24, 18).Value2 = "Hortaliza"

# --- New row 25 ---
$ws.Cells.Item(25, 1).Value2 = 9
$ws.Cells.Item(25, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(25, 3).Value2 = "Metropolitana"
$ws.Cells.Item(25, 4).Value2 = 44498
$ws.Cells.Item(25, 5).Value2 = 13
$ws.Cells.Item(25, 6).Value2 = 300000000
$ws.Cells.Item(25, 7).Value2 = "Espárragos"
$ws.Cells.Item(25, 8).Value2 = "Sin especificar"
$ws.Cells.Item(25, 9).Value2 = "Primera"
$ws.Cells.Item(25, 10).Value2 = 220
$ws.Cells.Item(25, 11).Value2 = 1000
$ws.Cells.Item(25, 12).Value2 = 1100
$ws.Cells.Item(25, 13).Value2 = 1055
$ws.Cells.Item(25, 14).Value2 = "`$/kilo"
$ws.Cells.Item(25, 15).Value2 = "Provincia de Linares"
$ws.Cells.Item(25, 16).Value2 = 1055
$ws.Cells.Item(25, 17).Value2 = 1
$ws.Cells.Item(25, 18).Value2 = "Hortaliza"

# --- New row 26 ---
$ws.Cells.Item(26, 1).Value2 = 9
$ws.Cells.Item(26, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(26, 3).Value2 = "Metropolitana"
$ws.Cells.Item(26, 4).Value2 = 44498
$ws.Cells.Item(26, 5).Value2 = 13
$ws.Cells.Item(26, 6).Value2 = 300000000
$ws.Cells.Item(26, 7).Value2 = "Espárragos"
$ws.Cells.Item(26, 8).Value2 = "Sin especificar"
$ws.Cells.Item(26, 9).Value2 = "Segunda"
$ws.Cells.Item(26, 10).Value2 = 120
$ws.Cells.Item(26, 11).Value2 = 800
$ws.Cells.Item(26, 12).Value2 = 900
$ws.Cells.Item(26, 13).Value2 = 858
$ws.Cells.Item(26, 14).Value2 = "`$/kilo"
$ws.Cells.Item(26, 15).Value2 = "Provincia de Linares"
$ws.Cells.Item(26, 16).Value2 = 858
$ws.Cells.Item(26, 17).Value2 = 1
$ws.Cells.Item(26, 18).Value2 = "Hortaliza"
